$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto quote data.
# Leading apostrophe forces Excel to store these numeric-looking strings as
# literal text (matching the inlineStr cells already in the sheet) instead of
# auto-converting them to numbers/percentages.
$ws.Range("D2").Value = "'290.79"
$ws.Range("E2").Value = "'-3.35%"
$ws.Range("D3").Value = "'30.72"
$ws.Range("E3").Value = "'-5.08%"
$ws.Range("D4").Value = "'4.954"
$ws.Range("E4").Value = "'0.03%"
$ws.Range("D5").Value = "'0.07219"
$ws.Range("E5").Value = "'-5.34%"
$ws.Range("D6").Value = "'1.863"
$ws.Range("E6").Value = "'-3.28%"
$ws.Range("D7").Value = "'7.692"
$ws.Range("E7").Value = "'-1.82%"
$ws.Range("D8").Value = "'3.771"
$ws.Range("E8").Value = "'-0.78%"
$ws.Range("D9").Value = "'0.8963"
$ws.Range("E9").Value = "'-2.38%"
$ws.Range("E10").Value = "'-5.05%"
$ws.Range("D11").Value = "'0.07699"
$ws.Range("E11").Value = "'-0.67%"
$ws.Range("D12").Value = "'0.07973"
$ws.Range("E12").Value = "'-6.58%"
$ws.Range("D13").Value = "'0.03035"
$ws.Range("E13").Value = "'-5.00%"
$ws.Range("D14").Value = "'0.1000"
$ws.Range("E14").Value = "'0.01%"
$ws.Range("D15").Value = "'0.001492"
$ws.Range("E15").Value = "'-1.18%"
$ws.Range("D16").Value = "'0.005729"
$ws.Range("E16").Value = "'-3.38%"
$ws.Range("E17").Value = "'-0.08%"
$ws.Range("D18").Value = "'3.468"
$ws.Range("E18").Value = "'0.12%"
$ws.Range("D19").Value = "'2.080"
$ws.Range("E19").Value = "'-3.31%"
$ws.Range("D20").Value = "'0.3319"
$ws.Range("D21").Value = "'0.1298"
$ws.Range("E21").Value = "'-2.16%"
$ws.Range("D22").Value = "'4.049"
$ws.Range("E22").Value = "'-5.24%"
$ws.Range("D23").Value = "'0.2321"
$ws.Range("E23").Value = "'16.51%"
$ws.Range("D24").Value = "'0.04513"
$ws.Range("E24").Value = "'-0.14%"
$ws.Range("D25").Value = "'0.001214"
$ws.Range("E25").Value = "'-0.66%"
$ws.Range("D26").Value = "'0.004640"
$ws.Range("E26").Value = "'5.75%"
$ws.Range("E27").Value = "'-0.13%"
$ws.Range("D39").Value = "'0.01579"
$ws.Range("E39").Value = "'-7.13%"
$ws.Range("E40").Value = "'-6.25%"
$ws.Range("D41").Value = "'0.007387"
$ws.Range("E41").Value = "'-1.22%"
$ws.Range("D42").Value = "'0.01003"
$ws.Range("E43").Value = "'-3.49%"
$ws.Range("D44").Value = "'0.002091"
$ws.Range("E44").Value = "'-10.42%"
$ws.Range("D45").Value = "'0.009190"
$ws.Range("E45").Value = "'-12.87%"
$ws.Range("D46").Value = "'0.00005964"
$ws.Range("E46").Value = "'-4.71%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("D48").Value = "'2.305"
$ws.Range("E48").Value = "'180.97%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.06%"
